# Scheduled-runner style refresh of computed price/profit columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# ARM, BSM, CUL and LTW sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# ARM sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H61").Value  = 4187960.5
$ws.Range("I61").Value  = 1812501
$ws.Range("K61").Value  = 1812501
$ws.Range("M61").Value  = -1812289

$ws.Range("H74").Value  = 26079934
$ws.Range("I74").Value  = 1021638.3
$ws.Range("J74").Value  = 76196520
$ws.Range("K74").Value  = 1021638.3
$ws.Range("L74").Value  = 76196520
$ws.Range("M74").Value  = -1020764.3
$ws.Range("N74").Value  = -76198268

$ws.Range("H77").Value  = 26079934
$ws.Range("I77").Value  = 1021638.3
$ws.Range("J77").Value  = 76196520
$ws.Range("K77").Value  = 5108191.5
$ws.Range("L77").Value  = 380982600
$ws.Range("M77").Value  = -5103823.5
$ws.Range("N77").Value  = -380991336

$ws.Range("H136").Value = 4187960.5
$ws.Range("I136").Value = 1812501
$ws.Range("K136").Value = 5437503
$ws.Range("M136").Value = -5434953

# ---------------------------------------------------------------
# BSM sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H105").Value = 2272.4736
$ws.Range("I105").Value = 2191.2144
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 2191.2144
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -444.2143999999998
$ws.Range("N105").Value = -5994

# These leves' price-lookup columns came back empty from the refresh;
# clear the stale cached numbers (rows 121 and 136 are unaffected).
$clearedRows = @(117,118,119,120,122,123,124,125,126,127,128,129,130,131,132,133,134,135,137,138,139,140,141)
foreach ($r in $clearedRows) {
    $rangeAddr = "H" + $r + ":N" + $r
    $ws.Range($rangeAddr).ClearContents()
}

# ---------------------------------------------------------------
# CUL sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H131").Value = 25889.5
$ws.Range("I131").Value = 111483.445
$ws.Range("J131").Value = 1039.6451
$ws.Range("K131").Value = 334450.335
$ws.Range("L131").Value = 3118.9353
$ws.Range("M131").Value = -329410.335
$ws.Range("N131").Value = -13198.9353

# ---------------------------------------------------------------
# LTW sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H132").Value = 2153688.2
$ws.Range("I132").Value = 2779117.8
$ws.Range("K132").Value = 8337353.399999999
$ws.Range("M132").Value = -8334823.399999999
